# Resets the value of items missing from the proposal.
# Column C ("Valor Prop.") has some cells that were left as an empty
# string (text) instead of a numeric 0. Excel stores these as shared
# string references; the fix replaces them with a literal numeric 0,
# keeping everything else (style, other columns) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column C is the 3rd column ("Valor Prop.") - data starts at row 2
# (row 1 is the header row).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()

    # Cells that are missing/blank proposal values were stored as an
    # empty string; reset them to numeric 0.
    if ($val -eq "") {
        $cell.Value = 0
    }
}
